$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -14.048
$ws.Range("C6").Value = -11.8327
$ws.Range("C7").Value = -11.91230000000001
$ws.Range("C16").Value = -11.55809999999999
$ws.Range("C20").Value = -14.5905
